$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.172.51"
$ws.Range("E2").Value = "  -1.00%  "

# Row 3
$ws.Range("D3").Value = "2.314.86"
$ws.Range("E3").Value = "  -1.95%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'312.74"
$ws.Range("E5").Value = "  -5.56%  "

# Row 6
$ws.Range("D6").Value = "'106.01"
$ws.Range("E6").Value = "  +6.15%  "

# Row 7
$ws.Range("E7").Value = "  -1.61%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  -2.63%  "

# Row 10
$ws.Range("D10").Value = "'40.12"
$ws.Range("E10").Value = "  +1.62%  "

# Row 11
$ws.Range("D11").Value = "'0.0916"
$ws.Range("E11").Value = "  -0.56%  "

# Row 12
$ws.Range("E12").Value = "  -1.68%  "

# Row 13
$ws.Range("E13").Value = "  +0.12%  "

# Row 14
$ws.Range("D14").Value = "'0.983"
$ws.Range("E14").Value = "  -2.30%  "

# Row 15
$ws.Range("D15").Value = "'15.58"
$ws.Range("E15").Value = "  -4.47%  "

# Row 16
$ws.Range("D16").Value = "2.660.80"
$ws.Range("E16").Value = "  -1.82%  "

# Row 17
$ws.Range("D17").Value = "2.324.08"
$ws.Range("E17").Value = "  -1.22%  "

# Row 18
$ws.Range("D18").Value = "42.154.71"
$ws.Range("E18").Value = "  -0.96%  "

# Row 19
$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "  -2.27%  "

# Row 20
$ws.Range("E20").Value = "  -1.49%  "

# Row 21
$ws.Range("E21").Value = "  -1.51%  "

# Row 22
$ws.Range("E22").Value = "  -6.91%  "

# Row 23
$ws.Range("D23").Value = "'259.07"
$ws.Range("E23").Value = "  -3.59%  "

# Row 24
$ws.Range("D24").Value = "'2.30"
$ws.Range("E24").Value = "  -0.33%  "

# Row 25
$ws.Range("D25").Value = "'9.26"
$ws.Range("E25").Value = "  -7.64%  "

# Row 26
$ws.Range("E26").Value = "  +0.48%  "

# Row 27
$ws.Range("D27").Value = "'11.02"
$ws.Range("E27").Value = "  -3.89%  "

# Row 28
$ws.Range("E28").Value = "  +3.29%  "

# Row 29
$ws.Range("D29").Value = "'22.85"
$ws.Range("E29").Value = "  -1.36%  "

# Row 30
$ws.Range("D30").Value = "'35.60"
$ws.Range("E30").Value = "  +0.83%  "

# Row 31
$ws.Range("D31").Value = "'0.0895"
$ws.Range("E31").Value = "  -0.62%  "

# Row 32
$ws.Range("D32").Value = "'163.44"
$ws.Range("E32").Value = "  -7.27%  "

# Row 33
$ws.Range("E33").Value = "  -5.45%  "

# Row 34
$ws.Range("E34").Value = "  -3.86%  "

# Row 35
$ws.Range("E35").Value = "  -2.53%  "

# Row 36
$ws.Range("E36").Value = "  +11.43%  "

# Row 37
$ws.Range("D37").Value = "'4.52"
$ws.Range("E37").Value = "  -1.68%  "

# Row 38
$ws.Range("D38").Value = "'0.0353"
$ws.Range("E38").Value = "  -1.42%  "

# Row 39
$ws.Range("D39").Value = "'2.77"
$ws.Range("E39").Value = "  -7.09%  "

# Row 40
$ws.Range("D40").Value = "'3.65"
$ws.Range("E40").Value = "  -4.01%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.48"
$ws.Range("E41").Value = "  -3.06%  "

# Row 42
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'70.65"
$ws.Range("E42").Value = "  +0.88%  "

# Row 43
$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "'96.91"
$ws.Range("E43").Value = "  +7.10%  "

# Row 44
$ws.Range("E44").Value = "  -1.73%  "

# Row 45
$ws.Range("E45").Value = "  +0.14%  "

# Row 46
$ws.Range("E46").Value = "  +2.71%  "

# Row 47
$ws.Range("E47").Value = "  -5.29%  "

# Row 48
$ws.Range("E48").Value = "  -1.59%  "

# Row 49
$ws.Range("D49").Value = "'8.95"
$ws.Range("E49").Value = "  -1.82%  "

# Row 50
$ws.Range("D50").Value = "'74.99"
$ws.Range("E50").Value = "  +6.99%  "

# Row 51
$ws.Range("E51").Value = "  -0.15%  "
